$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# header style already used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for columns I and J, rows 2..57
$iValues = @(7,6,7,7,6,7,7,7,7,8,7,6,5,7,6,8,6,6,9,6,7,7,6,7,8,8,7,5,7,6,5,8,7,6,6,8,6,8,8,7,6,7,6,5,7,5,6,4,6,8,8,4,6,8,1,5)
$jValues = @(7,7,7,7,7,8,7,7,7,9,7,6,5,8,6,8,7,6,9,7,7,7,6,7,8,8,8,6,7,7,6,8,8,6,7,8,6,8,8,8,7,8,6,6,8,7,6,5,7,8,8,4,6,8,1,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
